$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the curated dimension labels to measure labels
$ws.Range("C2").Value = "iaest-measure:rama"
$ws.Range("F2").Value = "iaest-measure:rama-descripcion"

# Row 3 "rama" / "rama-descripcion" columns now reference the measure type instead of dim
$ws.Range("C3").Value = "medida"
$ws.Range("F3").Value = "medida"

# Row 4 "rama" / "rama-descripcion" columns now reference xsd:int instead of skos:Concept
$ws.Range("C4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"

# Row 5 no longer carries the mapping file references for rama / rama-descripcion
$ws.Range("C5").Clear()
$ws.Range("F5").Clear()
